$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Content edits inside "Table1" (A1:C15)
#    Row 9 used to be the duplicate main()/int/"Insert the functions" row;
#    it becomes the new mainGame() row. Row 8 (ClearArea) gets re-cased to
#    clearArea(). color()/outputPosition() purposes get reworded.
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = "mainGame()"
$ws.Range("B9").Value = "void"
$ws.Range("C9").Value = "Insert the main program "
$ws.Range("A8").Value = "clearArea()"
$ws.Range("C10").Value = "Change output's color"
$ws.Range("C11").Value = "Set player's fixed position"

# ---------------------------------------------------------------------------
# 2) Banded-row fill colors + alignment across the table body
#    Header -> bold + dark orange fill, centered both ways
#    Even data rows -> dark orange fill (name/type plain, purpose centered)
#    Odd data rows  -> light orange fill (name/type plain, purpose centered)
# ---------------------------------------------------------------------------
$darkOrange  = 26316    # RGB(204,102,0)  -> FFCC6600
$lightOrange = 6737151  # RGB(255,204,102) -> FFFFCC66

$ws.Range("A1:C1").Interior.Color = $darkOrange
$ws.Range("A1:C1").HorizontalAlignment = -4108
$ws.Range("A1:C1").VerticalAlignment = -4108

$evenRows = @(2,4,6,8,10,12,14)
foreach ($r in $evenRows) {
    $ws.Range("A$r`:B$r").Interior.Color = $darkOrange
    $ws.Range("C$r").Interior.Color = $darkOrange
    $ws.Range("C$r").HorizontalAlignment = -4108
}

$oddRows = @(3,5,7,9,11,13,15)
foreach ($r in $oddRows) {
    $ws.Range("A$r`:B$r").Interior.Color = $lightOrange
    $ws.Range("C$r").Interior.Color = $lightOrange
    $ws.Range("C$r").HorizontalAlignment = -4108
}

# Row 11's Purpose cell also picked up vertical centering
$ws.Range("C11").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 3) Table style
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.TableStyle = "TableStyleMedium5"

# ---------------------------------------------------------------------------
# 4) Conditional formatting: 2-color scale over the whole table range
# ---------------------------------------------------------------------------
$rng = $ws.Range("A1:C15")
$cf = $rng.FormatConditions.AddColorScale(2)
$cf.ColorScaleCriteria(1).FormatColor.Color = 10285055  # FFFFEF9C
$cf.ColorScaleCriteria(2).FormatColor.Color = 8109667   # FF63BE7B

# ---------------------------------------------------------------------------
# 5) Selection moved to C12
# ---------------------------------------------------------------------------
$ws.Range("C12").Select()
